$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simple price/volume updates (row number matches sheet row)
$ws.Range("D2").Value = "28.476.83"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "1.866.69"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("E4").Value = "  -1.61%  "
$ws.Range("D5").Value = "'315.19"
$ws.Range("E5").Value = "  -0.96%  "
$ws.Range("E6").Value = "  -1.35%  "
$ws.Range("D7").Value = "'0.5071"
$ws.Range("E7").Value = "  -1.51%  "
$ws.Range("D8").Value = "'0.3901"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").Value = "'0.08342"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "'42.52"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").Value = "'1.103"
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("D12").Value = "'6.187"
$ws.Range("E12").Value = "  -1.21%  "
$ws.Range("D13").Value = "1.861.89"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").Value = "'20.28"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").Value = "'7.232"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "'1.009"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "'0.00001099"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "'91.26"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'0.06723"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "'17.63"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("E21").Value = "  -1.47%  "
$ws.Range("D22").Value = "'5.892"
$ws.Range("E22").Value = "  -1.42%  "
$ws.Range("D23").Value = "28.527.32"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'11.07"
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("D25").Value = "'2.201"
$ws.Range("D26").Value = "2.078.18"
$ws.Range("E26").Value = "  +2.74%  "
$ws.Range("D27").Value = "'157.33"
$ws.Range("E27").Value = "  -3.07%  "
$ws.Range("D28").Value = "'20.62"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("D29").Value = "'2.413"
$ws.Range("E29").Value = "  +1.95%  "
$ws.Range("D30").Value = "'126.36"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").Value = "'0.1037"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").Value = "'1.035"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "'5.739"
$ws.Range("D34").Value = "'3.620"
$ws.Range("E34").Value = "  -0.79%  "

# Rows 35/36 swapped: VeChain and Hedera exchanged positions, with new values
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.06635"
$ws.Range("E35").Value = "  +2.03%  "
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02449"
$ws.Range("E36").Value = "  +0.75%  "

$ws.Range("D37").Value = "'8.957"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").Value = "'0.2156"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").Value = "'5.023"
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").Value = "'1.178"
$ws.Range("E40").Value = "  -0.93%  "
$ws.Range("D41").Value = "'1.235"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").Value = "'0.6350"
$ws.Range("E42").Value = "  -1.39%  "
$ws.Range("D43").Value = "'11.07"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "'0.5991"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").Value = "'13.04"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("D47").Value = "'3.681"
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("D48").Value = "'1.995"
$ws.Range("E48").Value = "  +0.04%  "

# Rows 49/50 swapped: EOS and Quant exchanged positions, with new values
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'122.15"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.208"
$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("D51").Value = "'1.137"
$ws.Range("E51").Value = "  -7.42%  "
